$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.NumberFormat = "General"
}

$ws.Cells.Item(2, 4).Value = "42.901.91"
$ws.Cells.Item(2, 5).Value = "  +1.27%  "
$ws.Cells.Item(3, 4).Value = "2.288.83"
$ws.Cells.Item(3, 5).Value = "  +2.49%  "
$ws.Cells.Item(4, 5).Value = "  -0.02%  "
Set-TextCell 5 4 "252.25"
$ws.Cells.Item(5, 5).Value = "  +0.15%  "
Set-TextCell 6 4 "0.643"
$ws.Cells.Item(6, 5).Value = "  +3.35%  "
Set-TextCell 7 4 "73.74"
$ws.Cells.Item(7, 5).Value = "  +6.71%  "
Set-TextCell 9 4 "0.646"
$ws.Cells.Item(9, 5).Value = "  +2.58%  "
Set-TextCell 10 4 "39.23"
$ws.Cells.Item(10, 5).Value = "  -2.05%  "
Set-TextCell 11 4 "0.0980"
$ws.Cells.Item(11, 5).Value = "  +4.13%  "
Set-TextCell 12 4 "59.04"
$ws.Cells.Item(12, 5).Value = "  -0.68%  "
Set-TextCell 13 4 "7.43"
$ws.Cells.Item(13, 5).Value = "  +4.16%  "
Set-TextCell 14 4 "0.107"
$ws.Cells.Item(14, 5).Value = "  +1.92%  "
$ws.Cells.Item(15, 4).Value = "2.629.15"
$ws.Cells.Item(15, 5).Value = "  +2.20%  "
Set-TextCell 16 4 "15.32"
$ws.Cells.Item(16, 5).Value = "  +4.39%  "
Set-TextCell 17 4 "0.874"
$ws.Cells.Item(17, 5).Value = "  -0.64%  "
$ws.Cells.Item(18, 4).Value = "2.283.57"
$ws.Cells.Item(18, 5).Value = "  +2.17%  "
$ws.Cells.Item(19, 4).Value = "42.787.47"
$ws.Cells.Item(19, 5).Value = "  +1.13%  "
Set-TextCell 20 4 "0.0000101"
$ws.Cells.Item(20, 5).Value = "  +3.74%  "
Set-TextCell 21 4 "6.32"
$ws.Cells.Item(21, 5).Value = "  +1.93%  "
Set-TextCell 22 4 "72.59"
$ws.Cells.Item(22, 5).Value = "  -0.24%  "
Set-TextCell 23 4 "237.06"
$ws.Cells.Item(23, 5).Value = "  +2.23%  "
Set-TextCell 24 4 "2.23"
$ws.Cells.Item(24, 5).Value = "  +6.79%  "
$ws.Cells.Item(25, 5).Value = "  -1.26%  "
Set-TextCell 26 4 "11.55"
$ws.Cells.Item(26, 5).Value = "  +0.85%  "
$ws.Cells.Item(27, 5).Value = "  -0.07%  "
$ws.Cells.Item(28, 5).Value = "  -0.23%  "
$ws.Cells.Item(29, 5).Value = "  -0.69%  "
$ws.Cells.Item(30, 5).Value = "  -0.48%  "
Set-TextCell 31 4 "167.00"
$ws.Cells.Item(31, 5).Value = "  -0.33%  "
Set-TextCell 32 4 "21.05"
$ws.Cells.Item(32, 5).Value = "  +1.59%  "
Set-TextCell 33 4 "6.48"
$ws.Cells.Item(33, 5).Value = "  +6.15%  "
$ws.Cells.Item(34, 5).Value = "  +3.57%  "
Set-TextCell 35 4 "0.0825"
$ws.Cells.Item(35, 5).Value = "  +5.73%  "
Set-TextCell 36 4 "31.14"
$ws.Cells.Item(36, 5).Value = "  +11.73%  "
Set-TextCell 37 4 "0.126"
$ws.Cells.Item(37, 5).Value = "  +1.56%  "
Set-TextCell 38 4 "4.63"
$ws.Cells.Item(38, 5).Value = "  +12.64%  "
Set-TextCell 39 4 "4.77"
$ws.Cells.Item(39, 5).Value = "  +2.67%  "
$ws.Cells.Item(40, 5).Value = "  -3.26%  "
Set-TextCell 41 4 "14.42"
$ws.Cells.Item(41, 5).Value = "  +14.64%  "
Set-TextCell 42 4 "2.35"
$ws.Cells.Item(42, 5).Value = "  +3.91%  "
Set-TextCell 43 4 "5.94"
$ws.Cells.Item(43, 5).Value = "  +3.52%  "
Set-TextCell 44 4 "0.217"
$ws.Cells.Item(44, 5).Value = "  +8.99%  "
Set-TextCell 45 4 "61.89"
$ws.Cells.Item(45, 5).Value = "  -2.26%  "
Set-TextCell 46 4 "9.09"
$ws.Cells.Item(46, 5).Value = "  +4.55%  "
$ws.Cells.Item(47, 5).Value = "  -3.45%  "
$ws.Cells.Item(48, 5).Value = "  +2.06%  "
$ws.Cells.Item(49, 5).Value = "  +0.04%  "
$ws.Cells.Item(50, 2).Value = "Aave"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell 50 4 "101.21"
$ws.Cells.Item(50, 5).Value = "  +7.58%  "
$ws.Cells.Item(51, 2).Value = "ARBITRUM"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell 51 4 "1.18"
$ws.Cells.Item(51, 5).Value = "  +0.52%  "
